$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.968.65"
$ws.Range("E2").Value = "  +3.52%  "
$ws.Range("D3").Value = "2.614.72"
$ws.Range("E3").Value = "  +3.72%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.78"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.96"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "2.614.51"
$ws.Range("E9").Value = "  +3.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +13.78%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "3.098.22"
$ws.Range("E14").Value = "  +5.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.56"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("E16").Value = "  +6.73%  "
$ws.Range("D17").Value = "70.982.60"
$ws.Range("E17").Value = "  +4.19%  "
$ws.Range("D18").Value = "2.622.82"
$ws.Range("E18").Value = "  +5.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "380.67"
$ws.Range("E19").Value = "  +8.22%  "
$ws.Range("E20").Value = "  +4.43%  "
$ws.Range("E21").Value = "  +3.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.13"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("E24").Value = "  +4.76%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  +6.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("E27").Value = "  +6.27%  "
$ws.Range("E28").Value = "  +5.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +5.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "527.93"
$ws.Range("E31").Value = "  +4.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.01"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.31"
$ws.Range("E33").Value = "  +3.41%  "
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.92"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.119"
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("E39").Value = "  +6.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.95"
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("E43").Value = "  +8.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.02"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.330"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.12"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.99"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.63"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.531"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.67"
$ws.Range("E50").Value = "  +5.22%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0263"
$ws.Range("E51").Value = "  +1.01%  "
